$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I1").Value = "18/03/2023"

$ws.Range("C2").Value = 238
$ws.Range("D2").Value = 238
$ws.Range("E2").Value = 0
$ws.Range("F2").Value = 0
$ws.Range("G2").Value = 5
$ws.Range("I2").Value = 338
$ws.Range("J2").Value = -29.58579881656804

$ws.Range("C3").Value = 119
$ws.Range("D3").Value = 119
$ws.Range("F3").Value = 0
$ws.Range("G3").Value = 2
$ws.Range("I3").Value = 185
$ws.Range("J3").Value = -35.67567567567568

$ws.Range("C4").Value = 3
$ws.Range("D4").Value = 3
$ws.Range("E4").Value = 0
$ws.Range("I4").Value = 4
$ws.Range("J4").Value = -25

$ws.Range("B5").Value = 14
$ws.Range("C5").Value = 74
$ws.Range("D5").Value = 96
$ws.Range("E5").Value = 6
$ws.Range("F5").Value = 2
$ws.Range("I5").Value = 116
$ws.Range("J5").Value = -17.24137931034483

$ws.Range("C6").Value = 15
$ws.Range("D6").Value = 15
$ws.Range("I6").Value = 20
$ws.Range("J6").Value = -25

$ws.Range("B7").Value = 0
$ws.Range("C7").Value = 12
$ws.Range("D7").Value = 14
$ws.Range("E7").Value = 1
$ws.Range("G7").Value = 0
$ws.Range("I7").Value = 40
$ws.Range("J7").Value = -65

$ws.Range("C8").Value = 72
$ws.Range("D8").Value = 73
$ws.Range("E8").Value = 1
$ws.Range("F8").Value = 0
$ws.Range("G8").Value = 4
$ws.Range("I8").Value = 21
$ws.Range("J8").Value = 247.6190476190476

$ws.Range("C9").Value = 12
$ws.Range("D9").Value = 13
$ws.Range("E9").Value = 0
$ws.Range("I9").Value = 48
$ws.Range("J9").Value = -72.91666666666667

$ws.Range("B10").Value = 10
$ws.Range("C10").Value = 61
$ws.Range("D10").Value = 101
$ws.Range("E10").Value = 29
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 3
$ws.Range("I10").Value = 198
$ws.Range("J10").Value = -48.989898989899

$ws.Range("C11").Value = 91
$ws.Range("D11").Value = 92
$ws.Range("F11").Value = 0
$ws.Range("G11").Value = 1
$ws.Range("I11").Value = 127
$ws.Range("J11").Value = -27.55905511811023

$ws.Range("B12").Value = 2
$ws.Range("C12").Value = 131
$ws.Range("D12").Value = 183
$ws.Range("E12").Value = 13
$ws.Range("F12").Value = 0
$ws.Range("G12").Value = 2
$ws.Range("H12").Value = 38
$ws.Range("I12").Value = 256.1
$ws.Range("J12").Value = -28.54353768059352

$ws.Range("C13").Value = 6
$ws.Range("D13").Value = 7
$ws.Range("E13").Value = 0
$ws.Range("I13").Value = 263
$ws.Range("J13").Value = -97.33840304182509

$ws.Range("B14").Value = 2
$ws.Range("C14").Value = 143
$ws.Range("D14").Value = 178
$ws.Range("E14").Value = 11
$ws.Range("F14").Value = 0
$ws.Range("G14").Value = 3
$ws.Range("H14").Value = 23
$ws.Range("I14").Value = 344
$ws.Range("J14").Value = -48.25581395348837

$ws.Range("B15").Value = 1
$ws.Range("C15").Value = 59
$ws.Range("D15").Value = 63
$ws.Range("E15").Value = 2
$ws.Range("F15").Value = 1
$ws.Range("I15").Value = 102
$ws.Range("J15").Value = -38.23529411764706

$ws.Range("C16").Value = 60
$ws.Range("D16").Value = 78
$ws.Range("E16").Value = 18
$ws.Range("G16").Value = 0
$ws.Range("I16").Value = 92
$ws.Range("J16").Value = -15.21739130434783

$ws.Range("C17").Value = 14
$ws.Range("D17").Value = 17
$ws.Range("E17").Value = 3
$ws.Range("F17").Value = 0
$ws.Range("G17").Value = 0
$ws.Range("I17").Value = 57
$ws.Range("J17").Value = -70.17543859649122

$ws.Range("C18").Value = 0
$ws.Range("D18").Value = 0
$ws.Range("I18").Value = 2
$ws.Range("J18").Value = -100

$ws.Range("C19").Value = 4
$ws.Range("D19").Value = 4
$ws.Range("I19").Value = 5
$ws.Range("J19").Value = -20

$ws.Range("C20").Value = 15
$ws.Range("D20").Value = 17
$ws.Range("G20").Value = 0
$ws.Range("I20").Value = 22
$ws.Range("J20").Value = -22.72727272727273
